# "Close to ready for CRAN" - update the last benchmark row (row 26, httk 2.4.0)
# with corrected numbers after switching to chem props from ctxR, and leave a
# note about it in the Notes column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated benchmark statistics for the 2.4.0 row
$ws.Range("B26").Value = 1021

$ws.Range("F26").Value = 0.94769999999999999
$ws.Range("G26").Value = 353
$ws.Range("H26").Value = 0.27160000000000001
$ws.Range("I26").Value = 353

$ws.Range("J26").Value = 1.508
$ws.Range("K26").Value = 36

$ws.Range("L26").Value = 0.9698
$ws.Range("M26").Value = 80

$ws.Range("N26").Value = 1.1319999999999999
$ws.Range("O26").Value = 80

$ws.Range("P26").Value = 0.64659999999999995

# Leave a note explaining the change on this row
$ws.Range("R26").Value = "Switched to chem props from ctxR"

# Update the on-screen view/selection to match where the author left off
$ws.Range("F27").Select()

# Window was maximized before saving
$win = $wb.Windows.Item(1)
$win.Top = -110
$win.Left = -110
$win.Width = 19420
$win.Height = 11020
